$d = $word.ActiveDocument

# The "UI Sprites" list is one big paragraph; "UI_MAP_World_001" (World map)
# is now finished/implemented, so its row gets marked green, same as the
# other "done" rows elsewhere in the sheet (RGB 00B050). The following row
# ("UI_MAP_Lyndor_001") stays the default/automatic color.
#
# OLE_COLOR is 0x00BBGGRR, so RGB(0x00,0xB0,0x50) -> 0x0050B000.
$green = [System.Convert]::ToInt32("0050B000", 16)
$auto  = -16777216   # wdColorAutomatic

$text = $d.Content.Text
$worldStart   = $text.IndexOf("UI_MAP_World_001")
$lyndorStart  = $text.IndexOf("UI_MAP_Lyndor_001")
$creditsStart = $text.IndexOf("UI_CREDITS_concept_001")

# Whole "UI_MAP_World_001 ... 3hr" row (through its trailing line break).
$worldRange = $d.Range($worldStart, $lyndorStart)
$worldRange.Font.Color = $green

# Whole "UI_MAP_Lyndor_001 ... 3hr" row keeps the automatic color, but is
# now its own run once the row above has been split out and recolored.
$lyndorRange = $d.Range($lyndorStart, $creditsStart)
$lyndorRange.Font.Color = $auto
